# Generate Report for Handback
# - Flip the status text from "Ready for handoff" to "Handed back: in sync with en-US"
#   everywhere it appears (Overview summary + both language sheets).
# - Record a handback pass for each language sheet: populate the "Latest Target
#   File" / "Latest Handback File" columns (E/F) with the same file links as the
#   handoff columns (A/C), and stamp "Latest Handback DateTime" (G) with the
#   handback timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: the status column mirrors each language sheet's status ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

# ---- Per-language sheets ----
$langs = @(
    @{
        Name = "zh-cn";
        HandbackTime = "2016-03-10 07:36:06";
        Row2MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/fca907109b13182236542313bd145da6491b1de3/e2e/0ae8ddc0-ffb7-466d-8487-b778909afbaf.md";
        Row2XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b36e7afe9b1e5317d2244cf837f67dc8de2c8c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0ae8ddc0-ffb7-466d-8487-b778909afbaf.db711d57d5b216ad0a167c1feaf06132f22e064f.zh-cn.xlf";
        Row3MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/fca907109b13182236542313bd145da6491b1de3/e2e/34d9b026-1645-47bc-b32e-5a0e08b5037e.md";
        Row3XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b36e7afe9b1e5317d2244cf837f67dc8de2c8c3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/34d9b026-1645-47bc-b32e-5a0e08b5037e.68468a986e5d9ec7157df28f8072265916168e4f.zh-cn.xlf";
        Row2XlfName = "0ae8ddc0-ffb7-466d-8487-b778909afbaf.db711d57d5b216ad0a167c1feaf06132f22e064f.zh-cn.xlf";
        Row3XlfName = "34d9b026-1645-47bc-b32e-5a0e08b5037e.68468a986e5d9ec7157df28f8072265916168e4f.zh-cn.xlf";
    },
    @{
        Name = "de-de";
        HandbackTime = "2016-03-10 07:36:19";
        Row2MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/fca907109b13182236542313bd145da6491b1de3/e2e/0ae8ddc0-ffb7-466d-8487-b778909afbaf.md";
        Row2XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7d84f986d59c93dcb5d7baca8f791affb87ea70/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0ae8ddc0-ffb7-466d-8487-b778909afbaf.db711d57d5b216ad0a167c1feaf06132f22e064f.de-de.xlf";
        Row3MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/fca907109b13182236542313bd145da6491b1de3/e2e/34d9b026-1645-47bc-b32e-5a0e08b5037e.md";
        Row3XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7d84f986d59c93dcb5d7baca8f791affb87ea70/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/34d9b026-1645-47bc-b32e-5a0e08b5037e.68468a986e5d9ec7157df28f8072265916168e4f.de-de.xlf";
        Row2XlfName = "0ae8ddc0-ffb7-466d-8487-b778909afbaf.db711d57d5b216ad0a167c1feaf06132f22e064f.de-de.xlf";
        Row3XlfName = "34d9b026-1645-47bc-b32e-5a0e08b5037e.68468a986e5d9ec7157df28f8072265916168e4f.de-de.xlf";
    }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status -> handed back
    $ws.Range("B2").Value = $newStatus
    $ws.Range("B3").Value = $newStatus

    # Row 2 (0ae8ddc0-....md) — Latest Target File / Latest Handback File
    $ws.Hyperlinks.Add($ws.Range("E2"), $lang.Row2MdUrl, "", "", "0ae8ddc0-ffb7-466d-8487-b778909afbaf.md")
    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.Row2XlfUrl, "", "", $lang.Row2XlfName)

    # Row 3 (34d9b026-....md) — Latest Target File / Latest Handback File
    $ws.Hyperlinks.Add($ws.Range("E3"), $lang.Row3MdUrl, "", "", "34d9b026-1645-47bc-b32e-5a0e08b5037e.md")
    $ws.Hyperlinks.Add($ws.Range("F3"), $lang.Row3XlfUrl, "", "", $lang.Row3XlfName)

    # Latest Handback DateTime for both rows
    $ws.Range("G2").Value = $lang.HandbackTime
    $ws.Range("G3").Value = $lang.HandbackTime
}
